$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet keeps per-locale status in E2 (zh-cn) and F2 (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Locale detail sheets keep status in column C (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width shrink for the Status-related columns (narrower after text change) ---
# Overview: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5

Write-Host "Updated status text and column widths for localization report"
